$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths
$ws.Columns.Item(6).ColumnWidth = 3.140625
$ws.Columns.Item(8).ColumnWidth = 3.140625
$ws.Columns.Item(9).ColumnWidth = 2.140625
$ws.Columns.Item(10).ColumnWidth = 3.140625
$ws.Columns.Item(11).ColumnWidth = 5.7109375
$ws.Columns.Item(12).ColumnWidth = 5.7109375
$ws.Columns.Item(14).ColumnWidth = 5.7109375

# Update cell values in row 1
$ws.Range("A1").Value = 4
$ws.Range("B1").Value = 0
$ws.Range("C1").Value = 22
$ws.Range("D1").Value = 30
$ws.Range("E1").Value = 29
$ws.Range("F1").Value = 10
$ws.Range("G1").Value = 19
$ws.Range("H1").Value = 31
$ws.Range("I1").Value = 2
$ws.Range("J1").Value = 15
$ws.Range("K1").Value = 0.041999999999999996
$ws.Range("L1").Value = 0.096000000000000002
$ws.Range("M1").Value = 0.0019999999999999948
$ws.Range("N1").Value = 0.088999999999999996
$ws.Range("O1").Value = 0.025999999999999999
